# "moved some dates around" - update three due_date cells on Sheet1 and
# leave the selection where the author ended up (C24).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# C21: 2024-03-13 (45364) -> 2024-02-28 (45350)
$ws.Range("C21").Value = 45350

# C22: 2024-03-12 (45363) -> 2024-02-27 (45349)
$ws.Range("C22").Value = 45349

# C23: 2024-03-14 (45365) -> 2024-03-01 (45352)
# This cell also picks up an explicit black font color (distinct from the
# theme-based color used elsewhere), matching the author's edit.
$ws.Range("C23").Value = 45352
$ws.Range("C23").Font.Color = 0

# Leave the cursor on C24, matching the saved selection in the workbook.
$ws.Range("C24").Select()
